$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data point was inserted just before the old row 1033 entry,
# pushing the rest of the "Cilantro" series down by two rows (the series
# continues cyclically, so the two oldest rows reappear at the very end).
$ws.Range("A1033:A1034").EntireRow.Insert()

# Row 1033: new "Primera" / "$/caja 36 atados" reading
$ws.Range("A1033").Value = 9
$ws.Range("B1033").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1033").Value = "Metropolitana"
$ws.Range("D1033").Value = 45223
$ws.Range("E1033").Value = 13
$ws.Range("F1033").Value = 100112040
$ws.Range("G1033").Value = "Cilantro"
$ws.Range("H1033").Value = "Sin especificar"
$ws.Range("I1033").Value = "Primera"
$ws.Range("J1033").Value = 70
$ws.Range("K1033").Value = 8000
$ws.Range("L1033").Value = 8000
$ws.Range("M1033").Value = 8000
$ws.Range("N1033").Value = "$/caja 36 atados"
$ws.Range("O1033").Value = "Región Metropolitana"
$ws.Range("P1033").Value = 222
$ws.Range("Q1033").Value = 36
$ws.Range("R1033").Value = "Hortaliza"

# Row 1034: new "Primera" / "$/docena de atados" reading (same date)
$ws.Range("A1034").Value = 9
$ws.Range("B1034").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1034").Value = "Metropolitana"
$ws.Range("D1034").Value = 45223
$ws.Range("E1034").Value = 13
$ws.Range("F1034").Value = 100112040
$ws.Range("G1034").Value = "Cilantro"
$ws.Range("H1034").Value = "Sin especificar"
$ws.Range("I1034").Value = "Primera"
$ws.Range("J1034").Value = 160
$ws.Range("K1034").Value = 15000
$ws.Range("L1034").Value = 16000
$ws.Range("M1034").Value = 15500
$ws.Range("N1034").Value = "$/docena de atados"
$ws.Range("O1034").Value = "Región Metropolitana"
$ws.Range("P1034").Value = 5167
$ws.Range("Q1034").Value = 3
$ws.Range("R1034").Value = "Hortaliza"
